$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of trade data (row 5)
$ws.Range("A5").Value = 9920.67
$ws.Range("B5").Value = 9807.8799999999992
$ws.Range("C5").Value = 19.170000000000002
$ws.Range("D5").Value = 19.39
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = 1.1499999999999999
$ws.Range("G5").Value = 42609.505868055552
$ws.Range("H5").Value = $true

# Carry over the date/time number format from the row above (G column)
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
